$wb = $excel.ActiveWorkbook

# Update the "test" sheet's Course Color values (D1:F1) to RGB(128, 0, 128)
# Force text storage (matching original workbook's convention of storing
# these numeric-looking RGB components as text) via a temporary "@" number
# format, then clear the format so no style residue is left behind.
$testSheet = $wb.Worksheets.Item("test")
$testSheet.Range("D1:F1").NumberFormat = "@"
$testSheet.Range("D1").Value = "128"
$testSheet.Range("E1").Value = "0"
$testSheet.Range("F1").Value = "128"
$testSheet.Range("D1:F1").ClearFormats()

# Add a new worksheet named "fasf" after the "test" sheet
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add($null, $lastSheet)
$newSheet.Name = "fasf"
$newSheet.Range("A1").Value = "Secttion Number:"
$newSheet.Range("B1").Value = "dfs"
$newSheet.Range("C1").Value = "Course Color: "
$newSheet.Range("D1:F1").NumberFormat = "@"
$newSheet.Range("D1").Value = "65"
$newSheet.Range("E1").Value = "49"
$newSheet.Range("F1").Value = "140"
$newSheet.Range("D1:F1").ClearFormats()
